# Actualizacion automatica 2025-11-03 08:30:05
#
# This script re-runs the monthly rolling-window export:
#   - "VENTA MENSUAL" drops the oldest month (julio) and shifts
#     agosto/septiembre/octubre left by one column, adding a fresh
#     (still empty) "noviembre" column at the end.
#   - "VENTAS POR GRUPO" (the current-month-by-product-group sheet)
#     is rebuilt for the new month, so every cell that used to carry a
#     (now superseded) figure is reset to 0, and the "# de 35" summary
#     row is refreshed to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "VENTA MENSUAL": shift months C:F one column to the left.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Header row: each column's label moves to the previous column; a new
# month label ("noviembre") appears in column F.
$wsMensual.Range("C1").Value = $wsMensual.Range("D1").Value()
$wsMensual.Range("D1").Value = $wsMensual.Range("E1").Value()
$wsMensual.Range("E1").Value = $wsMensual.Range("F1").Value()
$wsMensual.Range("F1").Value = "noviembre"

# Data + totals rows (2..37): shift the four month columns left and
# zero out the newly opened month.
for ($r = 2; $r -le 37; $r++) {
    $d = $wsMensual.Cells.Item($r, 4).Value()
    $e = $wsMensual.Cells.Item($r, 5).Value()
    $f = $wsMensual.Cells.Item($r, 6).Value()
    $wsMensual.Cells.Item($r, 3).Value = $d
    $wsMensual.Cells.Item($r, 4).Value = $e
    $wsMensual.Cells.Item($r, 5).Value = $f
    $wsMensual.Cells.Item($r, 6).Value = 0
}

# The column widths were custom-fitted to the old figures; re-fit them
# to the new ones now occupying each column.
$wsMensual.Columns.Item(3).ColumnWidth = 13.166666666666666
$wsMensual.Columns.Item(4).ColumnWidth = 15.166666666666666
$wsMensual.Columns.Item(5).ColumnWidth = 13.166666666666666
$wsMensual.Columns.Item(6).ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------------
# 2) "VENTAS POR GRUPO": the current-month snapshot is regenerated, so
#    every product-group figure that belonged to the month that just
#    rolled off is cleared back to 0.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$cellsToClear = @(
    "L2", "M2",
    "D3", "K3", "M3",
    "L5", "M5",
    "D9", "L9", "M9",
    "M11",
    "D13", "M13",
    "M16",
    "H21", "I21", "L21", "M21", "P21",
    "M22",
    "L24",
    "M27",
    "E29", "I29",
    "M32",
    "M36"
)

foreach ($ref in $cellsToClear) {
    $wsGrupo.Range($ref).Value = 0
}

# Refresh the "# de 35" footer row (row 37) for every column touched
# above so the non-zero count reflects the cleared data.
$columnsToRefresh = @("D", "E", "H", "I", "K", "L", "M", "P")
foreach ($col in $columnsToRefresh) {
    $wsGrupo.Range($col + "37").Value = "0 de 35"
}
